$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the filename casing for the data importer configuration file
$ws.Range("B11").Value = "DataImportConfiguration.xlsx"

# Move the active selection to B11 (matches saved selection state in the diff)
$ws.Range("B11").Select()
